$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column (Price) cells are stored as literal text in the source data (even when
# they look numeric, e.g. "74.24"), so force text entry with a leading apostrophe and
# then clear the resulting quote-prefix formatting so the cell style is left untouched.
# E-column (Volume) cells are always non-numeric text (they carry padding spaces and a
# trailing "%"), so a plain .Value assignment is safe for those.

$ws.Range("D2").Value = "'42.190.93"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "'2.243.57"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'246.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("D7").Value = "'74.24"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.12%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.616"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.82%  "
$ws.Range("D10").Value = "'41.93"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.03%  "
$ws.Range("D11").Value = "'0.0945"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.61%  "
$ws.Range("D12").Value = "'7.10"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.28%  "
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("E14").Value = "  -3.28%  "
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "'2.225.64"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "'42.046.81"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "'71.94"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("E21").Value = "  +3.83%  "
$ws.Range("D22").Value = "'231.81"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("D23").Value = "'8.68"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +34.16%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'11.32"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("D26").Value = "'3.59"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.83%  "
$ws.Range("E27").Value = "  -2.83%  "
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").Value = "'169.33"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("D30").Value = "'20.62"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("E31").Value = "  -4.20%  "
$ws.Range("E32").Value = "  -3.46%  "
$ws.Range("D33").Value = "'30.50"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("D35").Value = "'5.16"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +9.24%  "
$ws.Range("D36").Value = "'4.51"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "'13.81"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").Value = "'2.19"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.26%  "
$ws.Range("D40").Value = "'5.78"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").Value = "'61.97"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").Value = "'0.202"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").Value = "'107.37"
$ws.Range("D43").ClearFormats()
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("E45").Value = "  -2.37%  "
$ws.Range("D46").Value = "'0.996"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").Value = "'4.32"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -8.82%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("E51").Value = "  +0.13%  "
